$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "XY Motor": update the axis-synchronisation calculation
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("XY Motor")

# --- updated input parameters -------------------------------------------------
$ws.Range("A7").Value = 1
$ws.Range("A13").Value = 64
$ws.Range("A15").Value = 200
$ws.Range("A17").Value = 200

# --- reworked "Count X" / "Count Y" formulas (Tabelle4) ----------------------
$ws.Range("B14").Formula = '=(C9*60*$A$9*10^6)/($A$17*D3*2*$A$13)'
$ws.Range("C14").Formula = '=(D9*60*$A$9*10^6)/($A$17*D4*2*$A$13)'

# --- new header row for the error-tracking table ------------------------------
$ws.Range("B16").Value = "Fehlerzeit X"
$ws.Range("C16").Value = "Fehlerzeit Y"
$ws.Range("D16").Value = "deltaX"
$ws.Range("E16").Value = "deltaY"
$ws.Range("F16").Value = "Fehler in Schritten"

# --- new data row with the error-tracking formulas ----------------------------
$ws.Range("B17").Formula = '=(Tabelle4[Round Count X]*A17*D3*2*A13)/(60*A9*10^6)'
$ws.Range("C17").Formula = '=(Tabelle4[Round Count Y]*A17*D4*2*A13)/(60*A9*10^6)'
$ws.Range("D17").Formula = '=C9-B17'
$ws.Range("E17").Formula = '=D9-C17'
$ws.Range("F17").Formula = '=D17/(1.6*10^-5)'

# --- turn B16:F17 into the new "Tabelle6" table -------------------------------
$lo6 = $ws.ListObjects.Add(1, $ws.Range("B16:F17"), 0, 1)
$lo6.Name = "Tabelle6"
$lo6.TableStyle = "TableStyleDark3"

# --- column width tweaks that came along with the new table ------------------
$ws.Range("B:C").ColumnWidth = 12.85546875
$ws.Range("F:F").ColumnWidth = 18.140625

# --- light-weight border accents mirroring the author's manual formatting ----
$ws.Range("A6").Borders.Item(8).LineStyle = 1
$ws.Range("A6").Borders.Item(7).LineStyle = 1
$ws.Range("A6").Borders.Item(10).LineStyle = 1

$ws.Range("A7:A17").Borders.Item(7).LineStyle = 1
$ws.Range("A17").Borders.Item(9).LineStyle = 1

$ws.Range("A16").Borders.Item(7).LineStyle = 1
$ws.Range("A16").Borders.Item(10).LineStyle = 1

$ws.Range("B3:D4").Borders.Item(7).LineStyle = 1
$ws.Range("E4").Borders.Item(7).LineStyle = 1

# -----------------------------------------------------------------
# Sheet "Auswahlregister": add the new prescaler option (64)
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Auswahlregister")
$ws2.Range("A4").Value = 64
$ws2.Range("A5").Select()

# -----------------------------------------------------------------
# Restore "XY Motor" as the active sheet/selection (as left by the author)
# -----------------------------------------------------------------
$ws.Activate()
$ws.Range("H18").Select()
